$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill_auto")

# Insert a new column before M (col 13) to make room for the new "regist" column.
$ws.Columns.Item(13).Insert()

# Fix up M3's style: Insert() copied L3's numeric style; match the blank style used
# by the other blank data cells in that column family (same as D3).
$ws.Range("D3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "L1=$($ws.Range('L1').Text)"
Write-Host "M1=$($ws.Range('M1').Text)"
Write-Host "N1=$($ws.Range('N1').Text)"
Write-Host "O1=$($ws.Range('O1').Text)"
Write-Host "dim=$($ws.UsedRange.Address())"
